# Daily attendance processing - 2025-11-18 11:21:03
# Reorders the "Recorded By" (column G) entries on the active sheet so that
# "System" is listed first, then "admin@admin.com", then the remaining
# recorder names in their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value2

    if ($value -eq $null -or $value -eq "") {
        continue
    }

    $parts = $value -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $indexed = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $token = $parts[$i]
        if ($token.Equals("System")) {
            $priority = 0
        } elseif ($token.Equals("admin@admin.com")) {
            $priority = 1
        } else {
            $priority = 2
        }
        $indexed += [PSCustomObject]@{ Priority = $priority; Order = $i; Token = $token }
    }

    # Combine priority + original order into a single sort key so the sort
    # is both correct and stable (ties broken by original position).
    $sorted = $indexed | Sort-Object -Property { $_.Priority * 1000 + $_.Order }
    $newValue = ($sorted | ForEach-Object { $_.Token }) -join ", "

    if ($newValue -ne $value) {
        $cell.Value2 = $newValue
    }
}
